# Fill in the Actual (F) and Correct (G) results for the
# "Fri, Jan 17, 2025" games (rows 99-109), and append the new
# "Sat, Jan 18, 2025" predictions (rows 110-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Results for games played Fri, Jan 17, 2025 (rows 99-109) ---
# Each entry: row, Actual (F), Correct (G)

$ws.Cells.Item(99,6).Value  = "Portland Winterhawks"
$ws.Cells.Item(99,7).Value  = 1

$ws.Cells.Item(100,6).Value = "Prince Albert Raiders"
$ws.Cells.Item(100,7).Value = 0

$ws.Cells.Item(101,6).Value = "Edmonton Oil Kings"
$ws.Cells.Item(101,7).Value = 1

$ws.Cells.Item(102,6).Value = "Saskatoon Blades"
$ws.Cells.Item(102,7).Value = 1

$ws.Cells.Item(103,6).Value = "Medicine Hat Tigers"
$ws.Cells.Item(103,7).Value = 0

$ws.Cells.Item(104,6).Value = "Vancouver Giants"
$ws.Cells.Item(104,7).Value = 1

$ws.Cells.Item(105,6).Value = "Lethbridge Hurricanes"
$ws.Cells.Item(105,7).Value = 1

$ws.Cells.Item(106,6).Value = "Everett Silvertips"
$ws.Cells.Item(106,7).Value = 0

$ws.Cells.Item(107,6).Value = "Seattle Thunderbirds"
$ws.Cells.Item(107,7).Value = 0

$ws.Cells.Item(108,6).Value = "Spokane Chiefs"
$ws.Cells.Item(108,7).Value = 1

$ws.Cells.Item(109,6).Value = "Victoria Royals"
$ws.Cells.Item(109,7).Value = 0

# --- New predictions for Sat, Jan 18, 2025 (rows 110-118) ---
# Columns: A GameID, B GameDate, C Home Team, D Away Team, E Prediction

$newDate = "Sat, Jan 18, 2025"

$ws.Cells.Item(110,1).Value = 1021658
$ws.Cells.Item(110,2).Value = $newDate
$ws.Cells.Item(110,3).Value = "Moose Jaw Warriors"
$ws.Cells.Item(110,4).Value = "Edmonton Oil Kings"
$ws.Cells.Item(110,5).Value = "Edmonton Oil Kings"

$ws.Cells.Item(111,1).Value = 1021659
$ws.Cells.Item(111,2).Value = $newDate
$ws.Cells.Item(111,3).Value = "Prince Albert Raiders"
$ws.Cells.Item(111,4).Value = "Red Deer Rebels"
$ws.Cells.Item(111,5).Value = "Prince Albert Raiders"

$ws.Cells.Item(112,1).Value = 1021662
$ws.Cells.Item(112,2).Value = $newDate
$ws.Cells.Item(112,3).Value = "Swift Current Broncos"
$ws.Cells.Item(112,4).Value = "Portland Winterhawks"
$ws.Cells.Item(112,5).Value = "Portland Winterhawks"

$ws.Cells.Item(113,1).Value = 1021657
$ws.Cells.Item(113,2).Value = $newDate
$ws.Cells.Item(113,3).Value = "Medicine Hat Tigers"
$ws.Cells.Item(113,4).Value = "Calgary Hitmen"
$ws.Cells.Item(113,5).Value = "Medicine Hat Tigers"

$ws.Cells.Item(114,1).Value = 1021656
$ws.Cells.Item(114,2).Value = $newDate
$ws.Cells.Item(114,3).Value = "Everett Silvertips"
$ws.Cells.Item(114,4).Value = "Prince George Cougars"
$ws.Cells.Item(114,5).Value = "Everett Silvertips"

$ws.Cells.Item(115,1).Value = 1021660
$ws.Cells.Item(115,2).Value = $newDate
$ws.Cells.Item(115,3).Value = "Seattle Thunderbirds"
$ws.Cells.Item(115,4).Value = "Kelowna Rockets"
$ws.Cells.Item(115,5).Value = "Seattle Thunderbirds"

$ws.Cells.Item(116,1).Value = 1021661
$ws.Cells.Item(116,2).Value = $newDate
$ws.Cells.Item(116,3).Value = "Spokane Chiefs"
$ws.Cells.Item(116,4).Value = "Lethbridge Hurricanes"
$ws.Cells.Item(116,5).Value = "Spokane Chiefs"

$ws.Cells.Item(117,1).Value = 1021663
$ws.Cells.Item(117,2).Value = $newDate
$ws.Cells.Item(117,3).Value = "Tri-City Americans"
$ws.Cells.Item(117,4).Value = "Wenatchee Wild"
$ws.Cells.Item(117,5).Value = "Wenatchee Wild"

$ws.Cells.Item(118,1).Value = 1021664
$ws.Cells.Item(118,2).Value = $newDate
$ws.Cells.Item(118,3).Value = "Victoria Royals"
$ws.Cells.Item(118,4).Value = "Kamloops Blazers"
$ws.Cells.Item(118,5).Value = "Victoria Royals"

# --- Update the view: scroll to show the new rows and select G118 ---
[void]$ws.Range("G118").Select()

Write-Host "Edit complete"
